$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Insert two new rows above the existing data row (current row 17) ---
# This pushes the existing "YAHYA" row down from 17 -> 19, and leaves a
# fresh blank row (copied formatting) at 17 and 18.
$ws.Rows.Item(17).Insert()
$ws.Rows.Item(17).Insert()

# Row 17 stays blank (already inherits the style of the row above via Insert).

# --- Row 18: brand new "TEST" room-registration entry ---
$ws.Range("A18").Value = "TEST"
$ws.Range("B18").Value = "TEST"
$ws.Range("C18").Value = "EA4C7814"
$ws.Range("D18").Value = 0
$ws.Range("E18").Value = "'100"
$ws.Range("F18").Value = "'TRUE"

# --- Row 19: existing "YAHYA" entry, now with an appended room and reset counter ---
$ws.Range("A19").Value = "YAHYA"
$ws.Range("B19").Value = "YAHYA"
$ws.Range("C19").Value = "EA4C7814"
$ws.Range("D19").Value = 0
$ws.Range("E19").Value = "'110"
$ws.Range("F19").Value = "'TRUE"

# The leading apostrophes above force Excel to store the numeric-looking
# values ("100", "110") and the boolean-looking value ("TRUE") as literal
# text (matching column C/E/F which are text columns elsewhere in the
# sheet), but they also mark the cells with a quote-prefix style. Re-apply
# the plain formatting from an already-clean cell in the same style family
# so the cells keep style "1" (left/wrap, no quote prefix) like the rest
# of the table.
$ws.Range("A14").Copy() | Out-Null
$ws.Range("E18").PasteSpecial(-4122) | Out-Null
$ws.Range("A14").Copy() | Out-Null
$ws.Range("F18").PasteSpecial(-4122) | Out-Null
$ws.Range("A14").Copy() | Out-Null
$ws.Range("E19").PasteSpecial(-4122) | Out-Null
$ws.Range("A14").Copy() | Out-Null
$ws.Range("F19").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false
